$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column N, mirroring the existing column M (2021 figures) with a
# new 2022 data point for each of the existing rows.

# Row 3: blank separator row, just needs the same (bottom-border) formatting.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# Row 4: header row with the year value.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2022

# Row 5: 2G coverage data.
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 98.8

# Row 6: 3G coverage data.
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 98

# Row 7: 4G coverage data.
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 96.9

$excel.CutCopyMode = $false

# Match the saved selection (active cell moved one column past the new data).
$null = $ws.Range("O4").Select()
